# update-mapping-pn13/ig/ValueSet-siph-typeeltpla-oncofair-valueset.xlsx
# "mise en commentaires des parametres de recherche qu'on n'utilise plus"
#
# Changes to the "Metadata" worksheet:
#   - the Date metadata value is refreshed to the new publication timestamp
#   - a new "Jurisdiction" property row is inserted right after "Contact"
#     (pushing Description / Purpose / Copyright / Immutable down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "Date" property value (row 8, column B) ---
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"

# --- Insert the new "Jurisdiction" row after "Contact" (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Give the new row the same look (font/fill/border/alignment) as the other
# data rows instead of the default formatting a freshly inserted row gets.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
